# Add "trapi_template" and "components" columns (with list data validations)
# to each of the *TestCase sheets, inserted right after "preconditions" (before
# the existing "id" column). This shifts id/name/description/tags from F:I to H:K.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "TestCase",
    "AcceptanceTestCase",
    "QuantitativeTestCase",
    "ComplianceTestCase",
    "KnowledgeGraphNavigationTestCase",
    "OneHopTestCase"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Insert two new blank columns at F and G; everything from the old F
    # column onward (id, name, description, tags) shifts right to H:K.
    $ws.Columns("F").Insert()
    $ws.Columns("G").Insert()

    $ws.Range("F1").Value = "trapi_template"
    $ws.Range("G1").Value = "components"

    $trapiRange = $ws.Range("F2:F1048576")
    $trapiRange.Validation.Add(3, 1, 1, '"ameliorates,treats,three_hop,drug_treats_rare_disease,drug-to-gene"')
    $trapiRange.Validation.ShowInput = $false
    $trapiRange.Validation.ShowError = $false

    $componentsRange = $ws.Range("G2:G1048576")
    $componentsRange.Validation.Add(3, 1, 1, '"arax,aragorn,ars,bte,improving"')
    $componentsRange.Validation.ShowInput = $false
    $componentsRange.Validation.ShowError = $false
}
